$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Dic_Disagg_Kategorien"

# New data set (KNr / Kategorie De / Kategorie En) replacing the old
# Einheit (ENr / Einheit De / Einheit En) rows.
$data = @(
    @("KNr", "Kategorie De", "Kategorie En"),
    @("K_AGE", "Altersklasse", "Age group"),
    @("K_AIRPOLL", "Luftschadstoff", "Air pollutant"),
    @("K_AREA", "Gebiet", "Area"),
    @("K_CALCMETH", "Berechnungsmethode", "Calculation method"),
    @("K_CRIM", "Straftat", "Criminal offence"),
    @("K_CRIMOFF", "Straftaten", "Criminal offences"),
    @("K_KREIS", "Kreis", "County"),
    @("K_LAENDER", "Bundesland", "Federal state"),
    @("K_PM", "Feinstaub", "Fine particulate matter"),
    @("K_SEA", "Meer", "Sea"),
    @("K_SERIES", "Zeitreihe", "Time series"),
    @("K_SEX", "Geschlecht", "Sex"),
    @("K_SUBINDEX", "Teilindizes", "Sub index"),
    @("K_TYPEAREA", "Art der Fläche", "Type of area"),
    @("K_URBAN", "Verstädterungsgrad", "Degree of urbanisation")
)

$oldLastRow = 28
$newLastRow = $data.Count

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Clear out the now-unused rows (old sheet had 28 rows, new one has 16)
if ($oldLastRow -gt $newLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 3))
    $clearRange.Clear()
}

# Adjust column widths to the new target widths (values chosen so the
# COM ColumnWidth -> stored OOXML width conversion lands on the target)
$ws.Columns.Item(1).ColumnWidth = 13.142857142857142
$ws.Columns.Item(2).ColumnWidth = 19.142857142857142
$ws.Columns.Item(3).ColumnWidth = 28.142857142857142
